# actualizacion 17 marzo con cambio de nps icx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumen")

function Set-TextValue($range, $text) {
    # Force the cell to store a literal text value (matching how the source
    # workbook was produced) instead of letting Excel auto-coerce a
    # numeric-looking string into a real number with an inferred format.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# ICX (row 17) - text values "4.80" -> "4.76"
Set-TextValue "B17" "4.76"
Set-TextValue "C17" "4.76"
Set-TextValue "D17" "4.76"

# NPS (row 18) - text values "81.96%" -> "85.71%"
Set-TextValue "B18" "85.71%"
Set-TextValue "C18" "85.71%"
Set-TextValue "D18" "85.71%"

# Numero de transacciones (row 20) - numeric values
$ws.Range("B20").Value = 10985752
$ws.Range("C20").Value = 13022348
$ws.Range("D20").Value = 24008100

# Valor transacciones (row 21) - text values with thousands separators
Set-TextValue "B21" "3,809,370,551,113"
Set-TextValue "C21" "4,388,130,730,022"
Set-TextValue "D21" "8,197,501,281,135"
